# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on all three sheets
#    (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2).
# 2) Narrow the "Status" column (Overview columns E/F, and column C on each
#    locale sheet) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status values -------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status columns --------------------------------------------
$newStatusWidth = 12.576851254417766

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
